$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw data path for the NVIS Extant intermediate file (B4):
# the source file was refreshed from a 20240709 export to a 20240801 export.
$ws.Range("B4").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_IUCNGET_DK_20240801.tif"

# Reselect the active cell to B4 to match the saved view state.
$ws.Range("B4").Select()
